$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (suchitra@gmail.com and suchitra1@gmail.com entries),
# shifting remaining rows up so that:
#   old A4 (suchitra2@gmail.com) -> A2
#   old A5 (suchitra4@gmail.com) -> A3
#   old A6 (duplicate roll number) -> A4
$ws.Rows("2:3").Delete()
